$d = $word.ActiveDocument
Write-Output "Content.End = $($d.Content.End)"
$r = $d.Range(200,200)
Write-Output "Range(200,200) start=$($r.Start) end=$($r.End)"
